$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E16:E20) is re-shuffled: the older periods 1710/1711/1712/1801/1802
# get reordered (newest first) while the matching "Valor Mora" (F) amounts follow the
# period they belong to. Only the period value, not its value, actually needs to move:
# period 1710 keeps its 8533 value (now on row 20) and the rest keep 32000 (now on row 16).
$ws.Range("E16").Value = "1802"
$ws.Range("E17").Value = "1801"
$ws.Range("E18").Value = "1712"
$ws.Range("E19").Value = "1711"
$ws.Range("E20").Value = "1710"

$ws.Range("F16").Value = 32000
$ws.Range("F17").Value = 32000
$ws.Range("F18").Value = 32000
$ws.Range("F19").Value = 32000
$ws.Range("F20").Value = 8533
